# Refresh the coin ranking table (rows 2-51) with the latest scrape results.
# Each row entry lists only the columns that changed for that coin; B/C (name/link)
# are plain text, D/E/G (price, 1h volume %, hour) are numeric-looking strings that
# must stay stored as text, matching the sheet's existing inlineStr cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row=2; D='317.39'; E='2.32%'; G='23' }
    @{ Row=3; D='41.05'; E='-0.13%'; G='23' }
    @{ Row=4; D='5.136'; E='0.29%'; G='23' }
    @{ Row=5; D='0.07638'; E='-0.62%'; G='23' }
    @{ Row=6; B='GateToken'; C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D='4.321'; E='0.67%'; G='23' }
    @{ Row=7; B='FTXToken'; C='https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D='1.685'; E='4.14%'; G='23' }
    @{ Row=8; B='MXToken'; C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D='0.9343'; E='1.45%'; G='23' }
    @{ Row=9; B='BTSEToken'; C='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; D='2.425'; E='-1.74%'; G='23' }
    @{ Row=10; B='LiechtensteinCryptoassetsExchange'; C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D='0.1244'; E='2.01%'; G='23' }
    @{ Row=11; B='WazirX'; C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D='0.1847'; E='1.31%'; G='23' }
    @{ Row=12; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.09063'; E='-1.09%'; G='23' }
    @{ Row=13; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.04116'; E='-4.81%'; G='23' }
    @{ Row=14; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.1057'; E='0.72%'; G='23' }
    @{ Row=15; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001292'; E='4.90%'; G='23' }
    @{ Row=16; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.005793'; E='-0.13%'; G='23' }
    @{ Row=17; B='UpBots'; C='https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'; D='0.007491'; E='1,897.31%'; G='23' }
    @{ Row=18; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.351'; E='-0.08%'; G='23' }
    @{ Row=19; D='0.3361'; E='1.45%'; G='23' }
    @{ Row=20; D='8.393'; E='21.62%'; G='23' }
    @{ Row=21; E='-2.96%'; G='23' }
    @{ Row=22; E='7.22%'; G='23' }
    @{ Row=23; D='0.04046'; E='0.35%'; G='23' }
    @{ Row=24; E='0.28%'; G='23' }
    @{ Row=25; D='0.004078'; E='-0.27%'; G='23' }
    @{ Row=26; E='0.49%'; G='23' }
    @{ Row=27; G='23' }
    @{ Row=28; G='23' }
    @{ Row=29; G='23' }
    @{ Row=30; G='23' }
    @{ Row=31; G='23' }
    @{ Row=32; G='23' }
    @{ Row=33; G='23' }
    @{ Row=34; G='23' }
    @{ Row=35; G='23' }
    @{ Row=36; G='23' }
    @{ Row=37; G='23' }
    @{ Row=38; D='0.02499'; E='1.43%'; G='23' }
    @{ Row=39; D='0.05244'; E='-0.31%'; G='23' }
    @{ Row=40; D='0.007760'; E='-0.85%'; G='23' }
    @{ Row=41; E='-1.13%'; G='23' }
    @{ Row=42; D='0.007057'; E='3.78%'; G='23' }
    @{ Row=43; D='0.002088'; E='13.38%'; G='23' }
    @{ Row=44; D='0.008226'; E='0.44%'; G='23' }
    @{ Row=45; D='0.3169'; E='2.08%'; G='23' }
    @{ Row=46; D='0.00006648'; E='-2.03%'; G='23' }
    @{ Row=47; E='0.58%'; G='23' }
    @{ Row=48; D='0.2337'; E='13.32%'; G='23' }
    @{ Row=49; E='3.03%'; G='23' }
    @{ Row=50; E='0.58%'; G='23' }
    @{ Row=51; E='0.58%'; G='23' }
)

$textForcedColumns = @('D', 'E', 'G')

foreach ($update in $rowUpdates) {
    $row = $update.Row
    foreach ($col in $update.Keys) {
        if ($col -eq 'Row') { continue }
        $cellRef = "$col$row"
        $newValue = $update[$col]
        if ($textForcedColumns -contains $col) {
            # Pre-format as Text so Excel doesn't coerce numeric-looking
            # strings (prices, percentages, hours) into real numbers,
            # then clear the formatting change so the cell style is left
            # exactly as it was (General / no explicit style).
            $ws.Range($cellRef).NumberFormat = '@'
            $ws.Range($cellRef).Value = $newValue
            $ws.Range($cellRef).ClearFormats()
        } else {
            $ws.Range($cellRef).Value = $newValue
        }
    }
}

